$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H1").Value = "coordxd"
$ws.Range("I1").Value = "coordyd"
$ws.Range("J1").Value = "p1"
$ws.Range("K1").Value = "p2"
$ws.Range("L1").Value = "p3"
$ws.Range("M1").Value = "p4"
$ws.Range("N1").Value = "Time"
$ws.Range("O1").Value = "Owner"

$ws.Range("O1").Select()
